$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- New row: Algorithm ---
$row1 = $t.Rows.Add()
$row1.Cells.Item(1).Range.Text = "Algorithm"

# The definition cell needs two separate runs ("...problems" + ".") so build
# it via a raw WordprocessingML fragment, then drop the leftover empty
# paragraph that InsertXML leaves behind ahead of our inserted paragraph.
$cell1b = $row1.Cells.Item(2)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>A finite sequence of well-defined instructions to solve a class of specific problems</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$cell1b.Range.InsertXML($xml1)
$cell1b.Range.Paragraphs.Item(1).Range.Delete()

# --- New row: Asymptotics ---
$row2 = $t.Rows.Add()
$row2.Cells.Item(1).Range.Text = "Asymptotics"
$row2.Cells.Item(2).Range.Text = "The growth of memory or program time when the problem size increases"

# The document used to end with two blank paragraphs after the table;
# the edit collapses that down to one. (Use Content.Paragraphs rather than
# Document.Paragraphs -- indexing the latter gets confused once a Tables
# reference has been touched in this session.)
$n = $d.Content.Paragraphs.Count
$d.Content.Paragraphs.Item($n - 1).Range.Delete()
